$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text: "[ug]" -> "[mg]"
$ws.Range("B1").Value = "Cysteine par semaine [mg]"

# Update nutrition values in column B (rows 2-79)
$ws.Range("B2").Value = 7353.945
$ws.Range("B3").Value = 7672.184999999999
$ws.Range("B4").Value = 4335.6
$ws.Range("B5").Value = 3875.4
$ws.Range("B6").Value = 7627.485
$ws.Range("B7").Value = 5349.235000000001
$ws.Range("B8").Value = 8963.235000000001
$ws.Range("B9").Value = 5923.885
$ws.Range("B10").Value = 9747.125
$ws.Range("B11").Value = 10637.2
$ws.Range("B12").Value = 6831.075000000001
$ws.Range("B13").Value = 5451.855
$ws.Range("B14").Value = 4455.12
$ws.Range("B15").Value = 3569.375
$ws.Range("B16").Value = 7652.69
$ws.Range("B17").Value = 10802.565
$ws.Range("B18").Value = 9005.385
$ws.Range("B19").Value = 14390.46
$ws.Range("B20").Value = 7458
$ws.Range("B21").Value = 3225.41
$ws.Range("B22").Value = 7169.825000000001
$ws.Range("B23").Value = 3110.8
$ws.Range("B24").Value = 8417.125
$ws.Range("B25").Value = 6211.790000000001
$ws.Range("B26").Value = 4319.775000000001
$ws.Range("B27").Value = 5083.535000000001
$ws.Range("B28").Value = 6189.835
$ws.Range("B29").Value = 7698
$ws.Range("B30").Value = 13668.565
$ws.Range("B31").Value = 7800.260000000001
$ws.Range("B32").Value = 3948.845
$ws.Range("B33").Value = 8522.385
$ws.Range("B34").Value = 4775.02
$ws.Range("B35").Value = 5541.02
$ws.Range("B36").Value = 3005.735
$ws.Range("B37").Value = 10315.885
$ws.Range("B38").Value = 20429.175
$ws.Range("B39").Value = 6587.050000000001
$ws.Range("B40").Value = 2597.6
$ws.Range("B41").Value = 6508.27
$ws.Range("B42").Value = 4887.55
$ws.Range("B43").Value = 6279.3
$ws.Range("B44").Value = 7156.674999999999
$ws.Range("B45").Value = 7843.895
$ws.Range("B46").Value = 6788.265
$ws.Range("B47").Value = 6176.884999999999
$ws.Range("B48").Value = 2197.85
$ws.Range("B49").Value = 9425.124999999998
$ws.Range("B50").Value = 3505.485
$ws.Range("B51").Value = 4885.235
$ws.Range("B52").Value = 4716.014999999999
$ws.Range("B53").Value = 3653.855
$ws.Range("B54").Value = 9308.309999999999
$ws.Range("B55").Value = 7796.135
$ws.Range("B56").Value = 4033.785
$ws.Range("B57").Value = 5510.035
$ws.Range("B58").Value = 10735.475
$ws.Range("B59").Value = 4730.785000000001
$ws.Range("B60").Value = 3919.185
$ws.Range("B61").Value = 8232.799999999999
$ws.Range("B62").Value = 6126.21
$ws.Range("B63").Value = 3424.42
$ws.Range("B64").Value = 2360.835
$ws.Range("B65").Value = 6612.315000000001
$ws.Range("B66").Value = 3970.7
$ws.Range("B67").Value = 6831.36
$ws.Range("B68").Value = 5915.915000000001
$ws.Range("B69").Value = 3191.02
$ws.Range("B70").Value = 2491.15
$ws.Range("B71").Value = 6689.285
$ws.Range("B72").Value = 6354.535000000001
$ws.Range("B73").Value = 2396.175
$ws.Range("B74").Value = 5733.56
$ws.Range("B75").Value = 5630.320000000001
$ws.Range("B76").Value = 5751.620000000001
$ws.Range("B77").Value = 8281.035
$ws.Range("B78").Value = 5348.474999999999
$ws.Range("B79").Value = 6368.110000000001
